# Updated CHE_grids model - 2025-08-19 17:28
# Re-shuffle the day/night timeslice groupings used by the EV charging
# use-case constraint, and the matching re_profiles hydro shares.

$wb = $excel.ActiveWorkbook

# --- ev_charging_uc sheet: update the day/night timeslice group lists ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")

# C13 ("D" group) and C14 ("N" group) hold comma separated timeslice lists
# that feed the G7 (=C14) and G8 (=C13) formulas below them.
$wsEv.Range("C13").Value = "RaD,WaD,FaD,FaP,SaP,RaP,SaD,WaP"
$wsEv.Range("C14").Value = "FaN,WaP,RaP,SaN,WaN,FaP,SaP,RaN"

# --- re_profiles sheet: re-order the hydro profile rows (M4:O7) ---
$wsRe = $wb.Worksheets.Item("re_profiles")

$wsRe.Range("M4").Value = "S"
$wsRe.Range("N4").Value = 0.39690767947648675

$wsRe.Range("M5").Value = "R"
$wsRe.Range("N5").Value = 0.27149547700006416

$wsRe.Range("M6").Value = "F"
$wsRe.Range("N6").Value = 0.27553730672996718

$wsRe.Range("M7").Value = "W"
$wsRe.Range("N7").Value = 0.2560595367934817

$wb.Save()
